# Updated cryptos list on Mon Aug 26 19:56:56 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (matches source inlineStr cells) instead of
# letting the numeric-looking strings (e.g. "552.69") get auto-coerced into
# a Number by Excel's usual type inference on Range.Value assignment.
function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "63.347.33"
$ws.Range("E2").Value = "  -1.21%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.683.70"
$ws.Range("E3").Value = "  -2.75%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.01%  "

# Row 5 - BNB
Set-TextValue "D5" "552.69"
$ws.Range("E5").Value = "  -4.10%  "

# Row 6 - Solana
Set-TextValue "D6" "158.09"
$ws.Range("E6").Value = "  -0.66%  "

# Row 7 - USDC
Set-TextValue "D7" "0.999"
$ws.Range("E7").Value = "  +0.07%  "

# Row 8 - XRP
Set-TextValue "D8" "0.583"
$ws.Range("E8").Value = "  -3.01%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -3.77%  "

# Row 10 - TRON
$ws.Range("E10").Value = "  -0.37%  "

# Row 12 - Toncoin
$ws.Range("E12").Value = "  -8.99%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "3.158.54"
$ws.Range("E13").Value = "  -2.78%  "

# Row 14 - Avalanche
Set-TextValue "D14" "26.30"
$ws.Range("E14").Value = "  -2.45%  "

# Row 15 - WrappedBTC
Set-TextValue "D15" "63.204.04"
$ws.Range("E15").Value = "  -0.85%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "2.687.65"
$ws.Range("E17").Value = "  -2.80%  "

# Row 18 - Chainlink
Set-TextValue "D18" "12.02"
$ws.Range("E18").Value = "  -0.88%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  -5.30%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "341.88"
$ws.Range("E20").Value = "  -4.83%  "

# Row 21 - Uniswap
Set-TextValue "D21" "6.31"
$ws.Range("E21").Value = "  -5.03%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  -0.40%  "

# Row 23 - Polygon
Set-TextValue "D23" "0.507"
$ws.Range("E23").Value = "  -4.07%  "

# Row 24 - Litecoin
Set-TextValue "D24" "63.82"
$ws.Range("E24").Value = "  -1.93%  "

# Row 25 - Kaspa
$ws.Range("E25").Value = "  -0.81%  "

# Row 26 - Binance-PegBSC-USD
$ws.Range("E26").Value = "  +0.08%  "

# Row 27 - InternetComputer(DFINITY)
Set-TextValue "D27" "8.16"
$ws.Range("E27").Value = "  -4.47%  "

# Row 28 - PEPE
Set-TextValue "D28" "0.0₃0853"
$ws.Range("E28").Value = "  -5.81%  "

# Row 29 - PancakeSwap
$ws.Range("E29").Value = "  -1.03%  "

# Row 30 - Fetch.AI
$ws.Range("E30").Value = "  -1.85%  "

# Row 31 - Aptos
$ws.Range("E31").Value = "  -5.09%  "

# Row 32 - Monero
Set-TextValue "D32" "166.15"
$ws.Range("E32").Value = "  -1.92%  "

# Row 34 - NEARProtocol
$ws.Range("E34").Value = "  -3.45%  "

# Row 35 - EthereumClassic
Set-TextValue "D35" "19.53"
$ws.Range("E35").Value = "  -3.38%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  -4.21%  "

# Row 37 - Stacks
$ws.Range("E37").Value = "  -2.61%  "

# Row 38 - Bittensor
Set-TextValue "D38" "339.25"
$ws.Range("E38").Value = "  -2.69%  "

# Row 39 - SuiNetwork
Set-TextValue "D39" "0.943"
$ws.Range("E39").Value = "  -6.11%  "

# Row 40 - RenderToken
Set-TextValue "D40" "6.06"
$ws.Range("E40").Value = "  -4.46%  "

# Row 41 - OKB
Set-TextValue "D41" "38.07"
$ws.Range("E41").Value = "  -2.69%  "

# Row 42 - Filecoin
$ws.Range("E42").Value = "  -6.21%  "

# Row 43 - InjectiveProtocol
Set-TextValue "D43" "20.76"
$ws.Range("E43").Value = "  -5.89%  "

# Row 44 - EnergySwap
$ws.Range("E44").Value = "  -6.21%  "

# Row 45 - Mantle
$ws.Range("E45").Value = "  -1.70%  "

# Row 46 - Hedera
$ws.Range("E46").Value = "  -4.53%  "

# Row 47 - FirstDigitalUSD
$ws.Range("E47").Value = "  +0.10%  "

# Row 48 - WhiteBITCoin
$ws.Range("E48").Value = "  +0.05%  "

# Row 49 - was Aave, now Stellar
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D49" "0.0971"
$ws.Range("E49").Value = "  -3.78%  "

# Row 50 - was Stellar, now Aave
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D50" "129.36"
$ws.Range("E50").Value = "  -6.09%  "

# Row 51 - Maker
Set-TextValue "D51" "2.094.69"
$ws.Range("E51").Value = "  -1.76%  "
